# Generate Report for Handback
#
# The 23a6492f-9465-453d-9c5d-4546e2c0cfe4.md record moves from
# "Ready for handoff" to "Handed back: in sync with en-US" and the three
# localization rows get re-sorted (23a6492f, ffff2d07819d, ffffffb95a6095)
# across the Overview / zh-cn / de-de sheets, picking up new handoff /
# handback file + datetime data along the way.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md"
$ws.Range("B2").Value = "e2e\23a6492f-9465-453d-9c5d-4546e2c0cfe4.md"
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"
$ws.Range("G2").Value = "2016-09-02 11:15:16"

$ws.Range("A3").Value = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"
$ws.Range("B3").Value = "e2e\ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-09-02 11:13:55"

$ws.Range("A4").Value = "ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md"
$ws.Range("B4").Value = "e2e\ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md"
$ws.Range("E4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "Handed back: in sync with en-US"
$ws.Range("G4").Value = "2016-09-02 11:13:55"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bef7e24114960b1d8375800025fc03010040d215/e2e/23a6492f-9465-453d-9c5d-4546e2c0cfe4.md", "", "", "e2e\23a6492f-9465-453d-9c5d-4546e2c0cfe4.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/63b77765f7eaf4f777f9f7472b1c961b7aa58ecc/e2e/ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md", "", "", "e2e\ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bef7e24114960b1d8375800025fc03010040d215/e2e/ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md", "", "", "e2e\ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 -> 23a6492f (newly handed back)
$ws.Range("A2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md"
$ws.Range("G2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.c51285a6a595635f217393391ba59a7ae22f88f2.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-02 11:15:00"
$ws.Range("I2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md"
$ws.Range("J2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.c51285a6a595635f217393391ba59a7ae22f88f2.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-02 11:15:53"

# Row 3 -> ffff2d07819d
$ws.Range("A3").Value = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"
$ws.Range("F3").Value = "False"

# Row 4 -> ffffffb95a6095
$ws.Range("A4").Value = "ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.zh-cn.xlf"
$ws.Range("H4").Value = "2016-09-02 11:13:51"
$ws.Range("I4").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.md"
$ws.Range("J4").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.zh-cn.xlf"
$ws.Range("K4").Value = "2016-09-02 11:14:17"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bef7e24114960b1d8375800025fc03010040d215/e2e/23a6492f-9465-453d-9c5d-4546e2c0cfe4.md", "", "", "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/58f3df30e9b561dcdecd75addcfc221ecc5f554a/e2e/23a6492f-9465-453d-9c5d-4546e2c0cfe4.md", "", "", "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/63b77765f7eaf4f777f9f7472b1c961b7aa58ecc/e2e/ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md", "", "", "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/58f3df30e9b561dcdecd75addcfc221ecc5f554a/e2e/49dbc38d-2320-4cc0-81ef-5fd882349c85.md", "", "", "49dbc38d-2320-4cc0-81ef-5fd882349c85.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bef7e24114960b1d8375800025fc03010040d215/e2e/ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md", "", "", "ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md")
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/58f3df30e9b561dcdecd75addcfc221ecc5f554a/e2e/49dbc38d-2320-4cc0-81ef-5fd882349c85.md", "", "", "49dbc38d-2320-4cc0-81ef-5fd882349c85.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 -> 23a6492f (newly handed back)
$ws.Range("A2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md"
$ws.Range("G2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.c51285a6a595635f217393391ba59a7ae22f88f2.de-de.xlf"
$ws.Range("H2").Value = "2016-09-02 11:15:16"
$ws.Range("I2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md"
$ws.Range("J2").Value = "23a6492f-9465-453d-9c5d-4546e2c0cfe4.c51285a6a595635f217393391ba59a7ae22f88f2.de-de.xlf"
$ws.Range("K2").Value = "2016-09-02 11:16:00"

# Row 3 -> ffff2d07819d
$ws.Range("A3").Value = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"

# Row 4 -> ffffffb95a6095
$ws.Range("A4").Value = "ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.de-de.xlf"
$ws.Range("H4").Value = "2016-09-02 11:13:55"
$ws.Range("I4").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.md"
$ws.Range("J4").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.de-de.xlf"
$ws.Range("K4").Value = "2016-09-02 11:14:25"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bef7e24114960b1d8375800025fc03010040d215/e2e/23a6492f-9465-453d-9c5d-4546e2c0cfe4.md", "", "", "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/592b0f96397ae94bfffca10150adf09b71a739e1/e2e/23a6492f-9465-453d-9c5d-4546e2c0cfe4.md", "", "", "23a6492f-9465-453d-9c5d-4546e2c0cfe4.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/63b77765f7eaf4f777f9f7472b1c961b7aa58ecc/e2e/ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md", "", "", "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/592b0f96397ae94bfffca10150adf09b71a739e1/e2e/49dbc38d-2320-4cc0-81ef-5fd882349c85.md", "", "", "49dbc38d-2320-4cc0-81ef-5fd882349c85.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bef7e24114960b1d8375800025fc03010040d215/e2e/ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md", "", "", "ffffffb95a6095-3218-4bfa-a2dc-c0713f6272ee.md")
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/592b0f96397ae94bfffca10150adf09b71a739e1/e2e/49dbc38d-2320-4cc0-81ef-5fd882349c85.md", "", "", "49dbc38d-2320-4cc0-81ef-5fd882349c85.md")

Write-Output "Generate Report for Handback: done"
